$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.578.24"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "2.645.55"
$ws.Range("E3").Value = "  -0.89%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.13"
$ws.Range("E5").Value = "  -0.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.97"

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.623"
$ws.Range("E8").Value = "  +3.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.126"
$ws.Range("E9").Value = "  +3.70%  "

$ws.Range("E10").Value = "  -0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  -2.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.156"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.61"
$ws.Range("E13").Value = "  -2.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000194"
$ws.Range("E14").Value = "  -0.22%  "

$ws.Range("D15").Value = "3.120.20"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").Value = "65.450.81"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "2.638.67"
$ws.Range("E17").Value = "  -1.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").Value = "  -1.53%  "

$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.47"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.87"
$ws.Range("E23").Value = "  -1.71%  "

$ws.Range("E24").Value = "  +3.00%  "

$ws.Range("E25").Value = "  -0.95%  "

$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.57"
$ws.Range("E27").Value = "  -2.74%  "

$ws.Range("E28").Value = "  -2.61%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.84"
$ws.Range("E30").Value = "  -2.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "529.63"
$ws.Range("E31").Value = "  -1.18%  "

$ws.Range("E32").Value = "  -1.34%  "

$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.44"
$ws.Range("E34").Value = "  +0.79%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.38"
$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.28"
$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.35"
$ws.Range("E40").Value = "  -3.83%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "159.92"
$ws.Range("E42").Value = "  -3.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.04"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.28"
$ws.Range("E44").Value = "  +1.41%  "

$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.46"
$ws.Range("E46").Value = "  -2.13%  "

$ws.Range("E47").Value = "  -2.36%  "

$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("E49").Value = "  +11.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0989"
$ws.Range("E50").Value = "  -0.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.55"
$ws.Range("E51").Value = "  -1.79%  "
